$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 715 (the "「ベイルート・アニメイテッド」" post), shifting
# all subsequent rows up by one.
$ws.Rows.Item(715).Delete()
